$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 5421.299537870102
$ws.Range("G4").Value = 3253.734504770488
$ws.Range("H4").Value = 1.666177596826546
$ws.Range("I4").Value = 9.686076806276105
$ws.Range("J4").Value = 58.23311956987891
$ws.Range("K4").Value = 3.371160993818194
$ws.Range("L4").Value = 75.68173810420558
$ws.Range("M4").Value = 24.33238917902054
$ws.Range("N4").Value = 14.82826949442824
$ws.Range("O4").Value = 6.616823517833836
$ws.Range("P4").Value = 278.0264926441014
